# Update column F (dSF) values for specific rows as per the data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -2
    9  = -3
    10 = -1
    25 = -7
    26 = 1
    33 = -3
    37 = 1
    38 = -7
    40 = -7
    43 = -7
    44 = 8
    47 = -8
    48 = 8
    53 = -2
    56 = -4
    57 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
